# Updated capital structure database
# Refresh computed ratio columns for the two Kuwait "Insurance (Prop/Cas.)"
# rows (row 2: country aggregate, row 3: First Takaful Insurance Company).
# Both rows receive the same recalculated set of values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2, 3) {
    $ws.Range("G$row").Value  = -0.0
    $ws.Range("H$row").Value  = -0.0
    $ws.Range("I$row").Value  = 1.103448275862069
    $ws.Range("J$row").Value  = 1.103448275862069
    $ws.Range("K$row").Value  = -15.9
    $ws.Range("L$row").Value  = 1.370689655172414

    $ws.Range("U$row").Value  = 0.02
    $ws.Range("V$row").Value  = 0.001538461538461538
    $ws.Range("W$row").Value  = -0.4746268656716418
    $ws.Range("X$row").Value  = 0.04466387908263446
    $ws.Range("Y$row").Value  = -0.5192907447542763
    $ws.Range("Z$row").Value  = -0.3464030817929346
    $ws.Range("AA$row").Value = -0.3822378833577209
    $ws.Range("AB$row").Value = 0.04466387908263446
    $ws.Range("AC$row").Value = -0.4269017624403554

    $ws.Range("AG$row").Value = -0.02
    $ws.Range("AJ$row").Value = -0.001540832049306626
    $ws.Range("AK$row").Value = -0.0006973500697350069
}

# O2/R2 go from 0 to -0 (sign only); O3/R3 go from -0 to 0.
$ws.Range("O2").Value = -0.0
$ws.Range("R2").Value = -0.0
$ws.Range("O3").Value = 0.0
$ws.Range("R3").Value = 0.0
